$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.792.45"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.108.70"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.97"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.26"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0779"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.417.00"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.58"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.16"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.130.66"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.676.89"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.23"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.08"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.63"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.03"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.45"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.58"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0623"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.50"
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.61"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.472.36"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.20"
$ws.Range("E46").Value = "  -10.22%  "
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.57"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.38"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.302.88"
$ws.Range("E51").Value = "  +2.08%  "
